# Dich edit_profile_screen.dart sang tieng Anh
# Appends a new "edit_profile_screen.dart" translation block (rows 32-45)
# to the existing EN/VI string table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: section header "edit_profile_screen.dart" (merged A32:B32, style like A1) ---
# Merge the (still-empty) destination first, THEN paste the format from the already-merged
# A1:B1 donor - pasting formats into cells before merging splits the donor's 4-side thin
# border into two half-borders (and forks new style/border entries), so order matters here.
$ws.Range("A32:B32").Merge() | Out-Null
$ws.Range("A1:B1").Copy() | Out-Null
$ws.Range("A32:B32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = "edit_profile_screen.dart"

# --- Row 33: gender options, stored with a leading apostrophe (quote-prefix / text-literal) ---
$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A33:B33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Formula = "'Male', 'Female', 'Other'"
$ws.Range("B33").Formula = "'Nam', 'Nữ', 'Khác'"

# --- Row 34: reuse the existing Save / Luu pair (style like A2) ---
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A34:B34").PasteSpecial(-4122) | Out-Null
$ws.Range("A34").Value = "Save"
$ws.Range("B34").Value = "Lưu"

# --- Rows 35-37: bold-ish sub header style (like row 28) ---
$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A35:B35").PasteSpecial(-4122) | Out-Null
$ws.Range("A35").Value = "Edit profile"
$ws.Range("B35").Value = "Chỉnh sửa thông tin"

$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A36:B36").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").Value = "Basic info "
$ws.Range("B36").Value = "Thông tin cơ bản"

$ws.Range("A28:B28").Copy() | Out-Null
$ws.Range("A37:B37").PasteSpecial(-4122) | Out-Null
$ws.Range("A37").Value = "Full name"
$ws.Range("B37").Value = "Họ và tên"

# --- Rows 38-45: plain rows (style like row 2) ---
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A38:B38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = "Gender"
$ws.Range("B38").Value = "Giới tính"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A39:B39").PasteSpecial(-4122) | Out-Null
$ws.Range("A39").Value = "Birthday"
$ws.Range("B39").Value = "Ngày sinh"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A40:B40").PasteSpecial(-4122) | Out-Null
$ws.Range("A40").Value = "Not selected"
$ws.Range("B40").Value = "Chưa có thông tin"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A41:B41").PasteSpecial(-4122) | Out-Null
$ws.Range("A41").Value = "Height"
$ws.Range("B41").Value = "Chiều cao"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A42:B42").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Value = "Weight"
$ws.Range("B42").Value = "Cân nặng"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A43:B43").PasteSpecial(-4122) | Out-Null
$ws.Range("A43").Value = "Interests & Style"
$ws.Range("B43").Value = "Sở thích & Phong cách"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A44:B44").PasteSpecial(-4122) | Out-Null
$ws.Range("A44").Value = "Personal style"
$ws.Range("B44").Value = "Phong cách cá nhân"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A45:B45").PasteSpecial(-4122) | Out-Null
$ws.Range("A45").Value = "Favorite colours"
$ws.Range("B45").Value = "Màu sắc yêu thích"

# --- View: move selection to the cell right after the new block ---
$ws.Range("A46").Select() | Out-Null

# --- Page setup: portrait A4, matching the rest of the workbook's print settings ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$excel.CutCopyMode = $false
